# ISIS1225 Lab 7 - "Entrega final"
# Adds a second machine's measurements ("Maquina 1" / "Maquina 2") below the
# existing PROBING / CHAINING tables, refreshes the already-collected data
# with the final measured values, and reproduces the light row-banding the
# author applied by hand to each 3-row data block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Duplicate the two existing blocks (header band + header row + 3
#    data rows) further down the sheet so the new machine's numbers
#    get exactly the same look & feel (fonts, number formats, etc.)
# ------------------------------------------------------------------
$ws.Range("A1:C5").Copy()
$ws.Range("A16").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A8:C12").Copy()
$ws.Range("A23").PasteSpecial(-4122)   # xlPasteFormats

$ws.Application.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Headers / titles for the new blocks (these all re-use existing
#    shared strings, mirroring the originals above).
# ------------------------------------------------------------------
$ws.Range("A16").Value = "Carga de Catálogo PROBING"
$ws.Range("A17").Value = "Factor de Carga (PROBING)"
$ws.Range("B17").Value = "Consumo de Datos [kB]"
$ws.Range("C17").Value = "Tiempo de Ejecución [ms]"

$ws.Range("A23").Value = "Carga de Catálogo CHAINING"
$ws.Range("A24").Value = "Factor de Carga (CHAINING)"
$ws.Range("B24").Value = "Consumo de Datos [kB]"
$ws.Range("C24").Value = "Tiempo de Ejecución [ms]"

$ws.Range("A16:C16").Merge()
$ws.Range("A23:C23").Merge()

# ------------------------------------------------------------------
# 3. Final measured values.
#    "Maquina 1" = the already-existing PROBING/CHAINING tables,
#    now refreshed with the final run's numbers.
#    "Maquina 2" = the new tables added further down the sheet.
# ------------------------------------------------------------------

# --- Maquina 1 : PROBING (rows 3-5) ---
$ws.Range("B3").Value = 466397.484
$ws.Range("C3").Value = 9531.2350000000006
$ws.Range("B4").Value = 466397.484
$ws.Range("C4").Value = 9578.5370000000003
$ws.Range("B5").Value = 466397.484
$ws.Range("C5").Value = 9747.8490000000002

# --- Maquina 1 : CHAINING (rows 10-12) ---
$ws.Range("B10").Value = 466397.484
$ws.Range("C10").Value = 9526.7340000000004
$ws.Range("B11").Value = 466397.484
$ws.Range("C11").Value = 9595.0040000000008
$ws.Range("B12").Value = 466397.484
$ws.Range("C12").Value = 9442.3780000000006

# --- Maquina 2 : PROBING (rows 18-20) ---
$ws.Range("A18").Value = 0.3
$ws.Range("B18").Value = 466353.29
$ws.Range("C18").Value = 9150.48
$ws.Range("A19").Value = 0.5
$ws.Range("B19").Value = 466353.29
$ws.Range("C19").Value = 9261.2009999999991
$ws.Range("A20").Value = 0.8
$ws.Range("B20").Value = 466353.29
$ws.Range("C20").Value = 9236.9709999999995

# --- Maquina 2 : CHAINING (rows 25-27) ---
$ws.Range("A25").Value = 2
$ws.Range("B25").Value = 466361.50099999999
$ws.Range("C25").Value = 9767.2990000000009
$ws.Range("A26").Value = 4
$ws.Range("B26").Value = 466361.50099999999
$ws.Range("C26").Value = 9585.1730000000007
$ws.Range("A27").Value = 6
$ws.Range("B27").Value = 466361.50099999999
$ws.Range("C27").Value = 9667.9449999999997

# ------------------------------------------------------------------
# 4. Machine labels next to each table's first data row.
# ------------------------------------------------------------------
$ws.Range("E3").Value = "Maquina 1"
$ws.Range("E18").Value = "Maquina 2"

# ------------------------------------------------------------------
# 5. Light banding on every 3-row data block (B:C columns): shaded /
#    plain / shaded-with-bottom-rule, General number format, explicit
#    black font (same font the header rows use). Every block is
#    formatted together in one shot so the four tables stay identical.
# ------------------------------------------------------------------
$bandAll   = $ws.Range("B3:C5,B10:C12,B18:C20,B25:C27")
$shaded    = $ws.Range("B3:C3,B10:C10,B18:C18,B25:C25")
$shadedBot = $ws.Range("B5:C5,B12:C12,B20:C20,B27:C27")

$bandAll.Style = "Normal"
$bandAll.Font.Color = 0
$bandAll.HorizontalAlignment = -4108
$bandAll.VerticalAlignment = -4108
$bandAll.WrapText = $true

$shaded.Interior.Color = 14277081
$shadedBot.Interior.Color = 14277081

$shadedBot.Borders.Item(9).Weight = -4138
$shadedBot.Borders.Item(9).Color = 0

# ------------------------------------------------------------------
# 6. Wrap the two new blocks in their own Excel Tables, matching the
#    look of the originals (Table1 / Table13).
# ------------------------------------------------------------------
$t3 = $ws.ListObjects.Add(1, $ws.Range("A17:C20"), $null, 1)
$t3.Name = "Table14"
$t3.TableStyle = "TableStyleLight1"

$t4 = $ws.ListObjects.Add(1, $ws.Range("A24:C27"), $null, 1)
$t4.Name = "Table135"
$t4.TableStyle = "TableStyleLight1"

# ------------------------------------------------------------------
# 7. Leave the selection where the author left it when saving.
# ------------------------------------------------------------------
$ws.Range("C30").Select()
